# Update the "Förändrad" date column (C) from 45179 (2023-09-10) to
# 45180 (2023-09-11) for rows 2 through 23, matching the automatic
# update reflected in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 23; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
